$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Add($ws.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee80dac869c4bd339375b9c6a2d69700c68e5c8c/e2e/6c4a31d1-9e8e-4c58-af82-47d40e29c6be.md", "", "", "6c4a31d1-9e8e-4c58-af82-47d40e29c6be.md")
